$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F14").Value = 1588
$ws.Range("F17").Value = 236
$ws.Range("F19").Value = 2044
$ws.Range("F20").Value = 1124
$ws.Range("F21").Value = 1019
$ws.Range("F24").Value = 815
$ws.Range("F25").Value = 1382
$ws.Range("F26").Value = 581
$ws.Range("F27").Value = 1314
$ws.Range("F29").Value = 272
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F13").Value = 26
$ws.Range("F25").Value = 82
$ws.Range("F37").Value = 49
$ws.Range("F39").Value = 279
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F17").Value = 324
$ws.Range("F18").Value = 161
$ws.Range("F19").Value = 617
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F10").Value = 324
$ws.Range("F13").Value = 161
$ws.Range("F15").Value = 1588
$ws.Range("F17").Value = 617
$ws.Range("F18").Value = 617
$ws.Range("F25").Value = 2044
$ws.Range("F26").Value = 1124
$ws.Range("F27").Value = 26
$ws.Range("F29").Value = 815
$ws.Range("F31").Value = 1382
$ws.Range("F35").Value = 581
$ws.Range("F39").Value = 1315
$ws.Range("F40").Value = 272
$ws.Range("F46").Value = 279
